$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.742.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.076.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "58.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0785"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.75%  "

$ws.Range("E11").Value = "  +2.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.383.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.060.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.675.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("E21").Value = "  +1.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.41%  "

$ws.Range("E27").Value = "  +4.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "

$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("E30").Value = "  -2.13%  "

$ws.Range("E31").Value = "  +2.60%  "

$ws.Range("E32").Value = "  +0.71%  "

$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.11%  "

$ws.Range("E35").Value = "  -4.58%  "

$ws.Range("E36").Value = "  +2.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.52%  "

$ws.Range("E38").Value = "  +0.09%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0969"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.33%  "

$ws.Range("E42").Value = "  -2.11%  "

$ws.Range("E43").Value = "  +0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.449.37"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.13%  "

$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.09%  "

$ws.Range("E48").Value = "  +0.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.95%  "

$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.267.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
